$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix the 14 swapped/rotated match rows (teams/odds/links corrected) ---
# Row 15
$ws.Cells.Item(15,6).Value = 'OFK Beograd'
$ws.Cells.Item(15,7).Value = 4
$ws.Cells.Item(15,8).Value = 'Metalac'
$ws.Cells.Item(15,10).Value = 1.87
$ws.Cells.Item(15,11).Value = '12/08/2023 23:08'
$ws.Cells.Item(15,12).Value = 1.83
$ws.Cells.Item(15,13).Value = '13/08/2023 17:15'
$ws.Cells.Item(15,14).Value = 3.15
$ws.Cells.Item(15,15).Value = '12/08/2023 23:08'
$ws.Cells.Item(15,16).Value = 3.62
$ws.Cells.Item(15,17).Value = '13/08/2023 17:15'
$ws.Cells.Item(15,18).Value = 3.85
$ws.Cells.Item(15,19).Value = '12/08/2023 23:08'
$ws.Cells.Item(15,20).Value = 3.52
$ws.Cells.Item(15,21).Value = '13/08/2023 17:15'
$ws.Cells.Item(15,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-metalac/rmtU5WuE/'

# Row 16
$ws.Cells.Item(16,6).Value = 'RFK Novi Sad'
$ws.Cells.Item(16,7).Value = 1
$ws.Cells.Item(16,8).Value = 'Radnicki Beograd'
$ws.Cells.Item(16,10).Value = 1.79
$ws.Cells.Item(16,11).Value = '12/08/2023 05:43'
$ws.Cells.Item(16,12).Value = 2.52
$ws.Cells.Item(16,13).Value = '13/08/2023 16:46'
$ws.Cells.Item(16,14).Value = 3.13
$ws.Cells.Item(16,15).Value = '12/08/2023 05:43'
$ws.Cells.Item(16,16).Value = 3.21
$ws.Cells.Item(16,17).Value = '13/08/2023 17:01'
$ws.Cells.Item(16,18).Value = 3.68
$ws.Cells.Item(16,19).Value = '12/08/2023 05:43'
$ws.Cells.Item(16,20).Value = 2.52
$ws.Cells.Item(16,21).Value = '13/08/2023 16:46'
$ws.Cells.Item(16,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-radnicki-beograd/SvbsrTm7/'

# Row 18
$ws.Cells.Item(18,6).Value = 'Metalac'
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = 'Tekstilac Odzaci'
$ws.Cells.Item(18,10).Value = 1.81
$ws.Cells.Item(18,11).Value = '19/08/2023 09:21'
$ws.Cells.Item(18,12).Value = 1.81
$ws.Cells.Item(18,13).Value = '19/08/2023 09:21'
$ws.Cells.Item(18,14).Value = 3.12
$ws.Cells.Item(18,15).Value = '19/08/2023 09:21'
$ws.Cells.Item(18,16).Value = 3.14
$ws.Cells.Item(18,17).Value = '19/08/2023 17:06'
$ws.Cells.Item(18,18).Value = 4.18
$ws.Cells.Item(18,19).Value = '19/08/2023 09:21'
$ws.Cells.Item(18,20).Value = 4.18
$ws.Cells.Item(18,21).Value = '19/08/2023 09:21'
$ws.Cells.Item(18,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/metalac-tekstilac-odzaci/ShwGdS20/'

# Row 19
$ws.Cells.Item(19,6).Value = 'Macva'
$ws.Cells.Item(19,7).Value = 3
$ws.Cells.Item(19,8).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(19,10).Value = 2.14
$ws.Cells.Item(19,11).Value = '18/08/2023 07:12'
$ws.Cells.Item(19,12).Value = 2.14
$ws.Cells.Item(19,13).Value = '18/08/2023 07:12'
$ws.Cells.Item(19,14).Value = 2.75
$ws.Cells.Item(19,15).Value = '18/08/2023 07:12'
$ws.Cells.Item(19,16).Value = 2.83
$ws.Cells.Item(19,17).Value = '19/08/2023 17:02'
$ws.Cells.Item(19,18).Value = 3.17
$ws.Cells.Item(19,19).Value = '18/08/2023 07:12'
$ws.Cells.Item(19,20).Value = 3.17
$ws.Cells.Item(19,21).Value = '18/08/2023 07:12'
$ws.Cells.Item(19,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-radnicki-s-mitrovica/YVwCc8mf/'

# Row 20
$ws.Cells.Item(20,6).Value = 'Sloboda'
$ws.Cells.Item(20,8).Value = 'Kolubara'
$ws.Cells.Item(20,10).Value = 2.58
$ws.Cells.Item(20,12).Value = 2.84
$ws.Cells.Item(20,13).Value = '19/08/2023 18:40'
$ws.Cells.Item(20,14).Value = 2.84
$ws.Cells.Item(20,16).Value = 2.84
$ws.Cells.Item(20,17).Value = '19/08/2023 18:40'
$ws.Cells.Item(20,18).Value = 2.68
$ws.Cells.Item(20,20).Value = 2.41
$ws.Cells.Item(20,21).Value = '19/08/2023 18:40'
$ws.Cells.Item(20,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/sloboda-kolubara/bw9EL1AP/'

# Row 25
$ws.Cells.Item(25,6).Value = 'FK Indjija'
$ws.Cells.Item(25,7).Value = 2
$ws.Cells.Item(25,8).Value = 'RFK Novi Sad'
$ws.Cells.Item(25,9).Value = 1
$ws.Cells.Item(25,10).Value = 1.45
$ws.Cells.Item(25,11).Value = '25/08/2023 08:13'
$ws.Cells.Item(25,12).Value = 1.34
$ws.Cells.Item(25,13).Value = '26/08/2023 16:45'
$ws.Cells.Item(25,14).Value = 3.74
$ws.Cells.Item(25,15).Value = '25/08/2023 08:13'
$ws.Cells.Item(25,16).Value = 4.28
$ws.Cells.Item(25,17).Value = '26/08/2023 16:45'
$ws.Cells.Item(25,18).Value = 5.1
$ws.Cells.Item(25,19).Value = '25/08/2023 08:13'
$ws.Cells.Item(25,20).Value = 8.119999999999999
$ws.Cells.Item(25,21).Value = '26/08/2023 16:45'
$ws.Cells.Item(25,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/indjija-rfk-novi-sad/OKWOf6XC/'

# Row 27
$ws.Cells.Item(27,6).Value = 'Smederevo'
$ws.Cells.Item(27,7).Value = 0
$ws.Cells.Item(27,8).Value = 'Macva'
$ws.Cells.Item(27,9).Value = 0
$ws.Cells.Item(27,10).Value = 1.81
$ws.Cells.Item(27,11).Value = '26/08/2023 13:43'
$ws.Cells.Item(27,12).Value = 2.25
$ws.Cells.Item(27,13).Value = '26/08/2023 16:55'
$ws.Cells.Item(27,14).Value = 3.15
$ws.Cells.Item(27,15).Value = '26/08/2023 13:43'
$ws.Cells.Item(27,16).Value = 2.92
$ws.Cells.Item(27,17).Value = '26/08/2023 16:55'
$ws.Cells.Item(27,18).Value = 3.97
$ws.Cells.Item(27,19).Value = '26/08/2023 13:43'
$ws.Cells.Item(27,20).Value = 3.14
$ws.Cells.Item(27,21).Value = '26/08/2023 16:55'
$ws.Cells.Item(27,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/smederevo-macva-sabac/l6zmaLOb/'

# Row 28
$ws.Cells.Item(28,6).Value = 'Mladost GAT'
$ws.Cells.Item(28,7).Value = 1
$ws.Cells.Item(28,8).Value = 'Vrsac'
$ws.Cells.Item(28,10).Value = 1.67
$ws.Cells.Item(28,12).Value = 1.78
$ws.Cells.Item(28,13).Value = '27/08/2023 16:51'
$ws.Cells.Item(28,14).Value = 3.27
$ws.Cells.Item(28,16).Value = 3.07
$ws.Cells.Item(28,17).Value = '27/08/2023 16:51'
$ws.Cells.Item(28,18).Value = 4.59
$ws.Cells.Item(28,20).Value = 4.62
$ws.Cells.Item(28,21).Value = '27/08/2023 16:51'
$ws.Cells.Item(28,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-vrsac/dWZjbuv5/'

# Row 30
$ws.Cells.Item(30,6).Value = 'OFK Beograd'
$ws.Cells.Item(30,7).Value = 2
$ws.Cells.Item(30,8).Value = 'Radnicki Beograd'
$ws.Cells.Item(30,10).Value = 1.68
$ws.Cells.Item(30,12).Value = 1.49
$ws.Cells.Item(30,13).Value = '27/08/2023 16:58'
$ws.Cells.Item(30,14).Value = 3.29
$ws.Cells.Item(30,16).Value = 3.84
$ws.Cells.Item(30,17).Value = '27/08/2023 16:58'
$ws.Cells.Item(30,18).Value = 4.48
$ws.Cells.Item(30,20).Value = 5.78
$ws.Cells.Item(30,21).Value = '27/08/2023 16:58'
$ws.Cells.Item(30,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-radnicki-beograd/Emrz2Nvt/'

# Row 31
$ws.Cells.Item(31,6).Value = 'Graficar Beograd'
$ws.Cells.Item(31,8).Value = 'Dubocica'
$ws.Cells.Item(31,10).Value = 1.69
$ws.Cells.Item(31,12).Value = 1.73
$ws.Cells.Item(31,13).Value = '27/08/2023 16:49'
$ws.Cells.Item(31,14).Value = 3.3
$ws.Cells.Item(31,16).Value = 3.29
$ws.Cells.Item(31,17).Value = '27/08/2023 16:49'
$ws.Cells.Item(31,18).Value = 4.38
$ws.Cells.Item(31,20).Value = 4.45
$ws.Cells.Item(31,21).Value = '27/08/2023 16:49'
$ws.Cells.Item(31,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-dubocica/2BVSgQnJ/'

# Row 33
$ws.Cells.Item(33,6).Value = 'Vrsac'
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(33,8).Value = 'FK Indjija'
$ws.Cells.Item(33,9).Value = 1
$ws.Cells.Item(33,10).Value = 2.27
$ws.Cells.Item(33,11).Value = '01/09/2023 05:13'
$ws.Cells.Item(33,12).Value = 2.71
$ws.Cells.Item(33,14).Value = 2.81
$ws.Cells.Item(33,15).Value = '01/09/2023 05:13'
$ws.Cells.Item(33,17).Value = '02/09/2023 16:03'
$ws.Cells.Item(33,18).Value = 2.86
$ws.Cells.Item(33,19).Value = '01/09/2023 05:13'
$ws.Cells.Item(33,20).Value = 2.61
$ws.Cells.Item(33,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/vrsac-indjija/zNYfcagB/'

# Row 34
$ws.Cells.Item(34,6).Value = 'Jedinstvo U.'
$ws.Cells.Item(34,7).Value = 3
$ws.Cells.Item(34,8).Value = 'Kolubara'
$ws.Cells.Item(34,9).Value = 2
$ws.Cells.Item(34,10).Value = 2.34
$ws.Cells.Item(34,11).Value = '02/09/2023 15:12'
$ws.Cells.Item(34,12).Value = 2.12
$ws.Cells.Item(34,14).Value = 2.84
$ws.Cells.Item(34,15).Value = '02/09/2023 15:12'
$ws.Cells.Item(34,17).Value = '02/09/2023 16:57'
$ws.Cells.Item(34,18).Value = 2.99
$ws.Cells.Item(34,19).Value = '02/09/2023 15:12'
$ws.Cells.Item(34,20).Value = 3.54
$ws.Cells.Item(34,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-kolubara/8CfpGuYn/'

# Row 42
$ws.Cells.Item(42,6).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(42,7).Value = 4
$ws.Cells.Item(42,8).Value = 'Jedinstvo U.'
$ws.Cells.Item(42,9).Value = 1
$ws.Cells.Item(42,10).Value = 2.04
$ws.Cells.Item(42,12).Value = 2.33
$ws.Cells.Item(42,13).Value = '09/09/2023 16:27'
$ws.Cells.Item(42,14).Value = 2.88
$ws.Cells.Item(42,16).Value = 3.03
$ws.Cells.Item(42,17).Value = '09/09/2023 15:46'
$ws.Cells.Item(42,18).Value = 3.23
$ws.Cells.Item(42,20).Value = 2.89
$ws.Cells.Item(42,21).Value = '09/09/2023 16:27'
$ws.Cells.Item(42,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-jedinstvo-ub/AumwAxwH/'

# Row 44
$ws.Cells.Item(44,6).Value = 'FK Indjija'
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = 'Graficar Beograd'
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 1.96
$ws.Cells.Item(44,12).Value = 2.09
$ws.Cells.Item(44,13).Value = '09/09/2023 13:14'
$ws.Cells.Item(44,14).Value = 3.12
$ws.Cells.Item(44,16).Value = 3.02
$ws.Cells.Item(44,17).Value = '09/09/2023 14:31'
$ws.Cells.Item(44,18).Value = 3.16
$ws.Cells.Item(44,20).Value = 3.33
$ws.Cells.Item(44,21).Value = '09/09/2023 13:14'
$ws.Cells.Item(44,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/indjija-graficar-beograd/xjeQCvgb/'

# --- Step 2: append 5 new matches (rows 60-64) played 22-23/09/2023, copying row 59's formatting ---
$ws.Range("A59:V59").Copy()
$ws.Range("A60:V64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 60
$ws.Cells.Item(60,1).Value = 59
$ws.Cells.Item(60,2).Value = 'serbia'
$ws.Cells.Item(60,3).Value = 'prva-liga'
$ws.Cells.Item(60,4).Value = '2023-2024'
$ws.Cells.Item(60,5).Value = 45192.66666666666
$ws.Cells.Item(60,6).Value = 'Mladost GAT'
$ws.Cells.Item(60,7).Value = 3
$ws.Cells.Item(60,8).Value = 'Jedinstvo U.'
$ws.Cells.Item(60,9).Value = 4
$ws.Cells.Item(60,10).Value = 1.82
$ws.Cells.Item(60,11).Value = '22/09/2023 03:13'
$ws.Cells.Item(60,12).Value = 1.9
$ws.Cells.Item(60,13).Value = '23/09/2023 05:06'
$ws.Cells.Item(60,14).Value = 2.99
$ws.Cells.Item(60,15).Value = '22/09/2023 03:13'
$ws.Cells.Item(60,16).Value = 3.07
$ws.Cells.Item(60,17).Value = '23/09/2023 15:38'
$ws.Cells.Item(60,18).Value = 3.78
$ws.Cells.Item(60,19).Value = '22/09/2023 03:13'
$ws.Cells.Item(60,20).Value = 3.93
$ws.Cells.Item(60,21).Value = '23/09/2023 15:38'
$ws.Cells.Item(60,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-jedinstvo-ub/vgbXNPId/'

# Row 61
$ws.Cells.Item(61,1).Value = 60
$ws.Cells.Item(61,2).Value = 'serbia'
$ws.Cells.Item(61,3).Value = 'prva-liga'
$ws.Cells.Item(61,4).Value = '2023-2024'
$ws.Cells.Item(61,5).Value = 45192.66666666666
$ws.Cells.Item(61,6).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(61,7).Value = 1
$ws.Cells.Item(61,8).Value = 'Dubocica'
$ws.Cells.Item(61,9).Value = 1
$ws.Cells.Item(61,10).Value = 1.85
$ws.Cells.Item(61,11).Value = '22/09/2023 03:13'
$ws.Cells.Item(61,12).Value = 1.64
$ws.Cells.Item(61,13).Value = '23/09/2023 15:55'
$ws.Cells.Item(61,14).Value = 2.98
$ws.Cells.Item(61,15).Value = '22/09/2023 03:13'
$ws.Cells.Item(61,16).Value = 3.78
$ws.Cells.Item(61,17).Value = '23/09/2023 15:55'
$ws.Cells.Item(61,18).Value = 3.66
$ws.Cells.Item(61,19).Value = '22/09/2023 03:13'
$ws.Cells.Item(61,20).Value = 4.33
$ws.Cells.Item(61,21).Value = '23/09/2023 15:55'
$ws.Cells.Item(61,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-dubocica/dUYnGDaj/'

# Row 62
$ws.Cells.Item(62,1).Value = 61
$ws.Cells.Item(62,2).Value = 'serbia'
$ws.Cells.Item(62,3).Value = 'prva-liga'
$ws.Cells.Item(62,4).Value = '2023-2024'
$ws.Cells.Item(62,5).Value = 45192.66666666666
$ws.Cells.Item(62,6).Value = 'Smederevo'
$ws.Cells.Item(62,7).Value = 2
$ws.Cells.Item(62,8).Value = 'Radnicki Beograd'
$ws.Cells.Item(62,9).Value = 1
$ws.Cells.Item(62,10).Value = 1.71
$ws.Cells.Item(62,11).Value = '22/09/2023 03:13'
$ws.Cells.Item(62,12).Value = 1.78
$ws.Cells.Item(62,13).Value = '23/09/2023 15:38'
$ws.Cells.Item(62,14).Value = 3.17
$ws.Cells.Item(62,15).Value = '22/09/2023 03:13'
$ws.Cells.Item(62,16).Value = 3.28
$ws.Cells.Item(62,17).Value = '23/09/2023 15:38'
$ws.Cells.Item(62,18).Value = 4.04
$ws.Cells.Item(62,19).Value = '22/09/2023 03:13'
$ws.Cells.Item(62,20).Value = 4.22
$ws.Cells.Item(62,21).Value = '23/09/2023 15:38'
$ws.Cells.Item(62,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/smederevo-radnicki-beograd/zLXjFXEd/'

# Row 63
$ws.Cells.Item(63,1).Value = 62
$ws.Cells.Item(63,2).Value = 'serbia'
$ws.Cells.Item(63,3).Value = 'prva-liga'
$ws.Cells.Item(63,4).Value = '2023-2024'
$ws.Cells.Item(63,5).Value = 45192.66666666666
$ws.Cells.Item(63,6).Value = 'Tekstilac Odzaci'
$ws.Cells.Item(63,7).Value = 4
$ws.Cells.Item(63,8).Value = 'Graficar Beograd'
$ws.Cells.Item(63,9).Value = 0
$ws.Cells.Item(63,10).Value = 2.16
$ws.Cells.Item(63,11).Value = '22/09/2023 03:13'
$ws.Cells.Item(63,12).Value = 2.58
$ws.Cells.Item(63,13).Value = '23/09/2023 15:48'
$ws.Cells.Item(63,14).Value = 2.99
$ws.Cells.Item(63,15).Value = '22/09/2023 03:13'
$ws.Cells.Item(63,16).Value = 2.88
$ws.Cells.Item(63,17).Value = '23/09/2023 15:48'
$ws.Cells.Item(63,18).Value = 2.86
$ws.Cells.Item(63,19).Value = '22/09/2023 03:13'
$ws.Cells.Item(63,20).Value = 2.71
$ws.Cells.Item(63,21).Value = '23/09/2023 15:48'
$ws.Cells.Item(63,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-graficar-beograd/SfUrHgpp/'

# Row 64
$ws.Cells.Item(64,1).Value = 63
$ws.Cells.Item(64,2).Value = 'serbia'
$ws.Cells.Item(64,3).Value = 'prva-liga'
$ws.Cells.Item(64,4).Value = '2023-2024'
$ws.Cells.Item(64,5).Value = 45192.66666666666
$ws.Cells.Item(64,6).Value = 'Vrsac'
$ws.Cells.Item(64,7).Value = 1
$ws.Cells.Item(64,8).Value = 'Metalac'
$ws.Cells.Item(64,9).Value = 0
$ws.Cells.Item(64,10).Value = 2.12
$ws.Cells.Item(64,11).Value = '22/09/2023 04:12'
$ws.Cells.Item(64,12).Value = 2.23
$ws.Cells.Item(64,13).Value = '23/09/2023 08:37'
$ws.Cells.Item(64,14).Value = 2.76
$ws.Cells.Item(64,15).Value = '22/09/2023 04:12'
$ws.Cells.Item(64,16).Value = 2.81
$ws.Cells.Item(64,17).Value = '23/09/2023 15:03'
$ws.Cells.Item(64,18).Value = 3.21
$ws.Cells.Item(64,19).Value = '22/09/2023 04:12'
$ws.Cells.Item(64,20).Value = 3.29
$ws.Cells.Item(64,21).Value = '23/09/2023 08:37'
$ws.Cells.Item(64,22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/vrsac-metalac/jT1yNqY2/'
